$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.211.37'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '1.558.02'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = '''288.22'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '''0.3812'
$ws.Range('E7').Value = '  +2.28%  '
$ws.Range('D8').Value = '''0.3321'
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').Value = '''44.75'
$ws.Range('D10').Value = '''1.141'
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').Value = '''0.07406'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('D12').Value = '''1.002'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '''20.20'
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('D14').Value = '''5.842'
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').Value = '''6.746'
$ws.Range('E15').Value = '  -2.43%  '
$ws.Range('D16').Value = '1.563.94'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').Value = '''0.00001074'
$ws.Range('E17').Value = '  -3.78%  '
$ws.Range('D18').Value = '''0.06652'
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('D19').Value = '''86.46'
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('D20').Value = '''6.403'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').Value = '''16.14'
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('D23').Value = '''11.73'
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('D24').Value = '22.200.39'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').Value = '''2.299'
$ws.Range('E25').Value = '  -3.71%  '
$ws.Range('D26').Value = '''2.556'
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').Value = '''151.20'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').Value = '''19.23'
$ws.Range('E28').Value = '  -2.36%  '
$ws.Range('D29').Value = '''4.942'
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('D30').Value = '''123.15'
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('D31').Value = '1.737.33'
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').Value = '''1.090'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('D33').Value = '''5.915'
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('D34').Value = '''1.918'
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('D35').Value = '''0.08215'
$ws.Range('E35').Value = '  -1.06%  '
$ws.Range('D36').Value = '''9.308'
$ws.Range('E36').Value = '  -3.79%  '
$ws.Range('D37').Value = '''0.06325'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').Value = '''0.02329'
$ws.Range('E38').Value = '  -5.21%  '
$ws.Range('D39').Value = '''5.317'
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').Value = '''0.2162'
$ws.Range('E40').Value = '  -4.90%  '
$ws.Range('D41').Value = '''1.232'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('D42').Value = '''10.98'
$ws.Range('E42').Value = '  -2.77%  '
$ws.Range('D43').Value = '''0.6068'
$ws.Range('D44').Value = '''1.001'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '''13.75'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').Value = '''3.744'
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('D47').Value = '''0.5866'
$ws.Range('E47').Value = '  -4.56%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.964'
$ws.Range('E48').Value = '  -4.13%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''122.01'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').Value = '''1.178'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('D51').Value = '''0.07060'
$ws.Range('E51').Value = '  -2.88%  '
